# Auto-generated Excel COM-interop script
# Rebrand-related data refresh: scenario_comparison.xlsx values updated
$wb = $excel.ActiveWorkbook

$wsRaw = $wb.Worksheets.Item("Raw_Data")
$wsSummary = $wb.Worksheets.Item("Summary")

# --- Raw_Data sheet updates ---
$wsRaw.Range("B2").Value = 0.8678757465277042
$wsRaw.Range("C2").Value = 0.5149789247094966
$wsRaw.Range("D2").Value = 0.5042681511465438
$wsRaw.Range("G2").Value = 0.06916539691385787
$wsRaw.Range("H2").Value = 5.13538877421817
$wsRaw.Range("L2").Value = 112.3458049393424
$wsRaw.Range("A3").Value = 593
$wsRaw.Range("B3").Value = 30.71866907403464
$wsRaw.Range("C3").Value = 3.081840487483596
$wsRaw.Range("D3").Value = 3.081804573429415
$wsRaw.Range("E3").Value = 24.18330223700058
$wsRaw.Range("F3").Value = 24.10613305948958
$wsRaw.Range("G3").Value = 2.383401061637665
$wsRaw.Range("H3").Value = 89.65922784542758
$wsRaw.Range("J3").Value = 1.867900580113075
$wsRaw.Range("L3").Value = 450.0812534530714
$wsRaw.Range("A4").Value = 588
$wsRaw.Range("B4").Value = 71.8903008454833
$wsRaw.Range("C4").Value = 4.734594539096623
$wsRaw.Range("D4").Value = 4.734099068704251
$wsRaw.Range("E4").Value = 38.53131181010473
$wsRaw.Range("F4").Value = 38.45428315748122
$wsRaw.Range("G4").Value = 3.818851138070488
$wsRaw.Range("H4").Value = 177.0949628351989
$wsRaw.Range("J4").Value = 3.689478392399978
$wsRaw.Range("L4").Value = 1376.955136523368
$wsRaw.Range("A5").Value = 575
$wsRaw.Range("B5").Value = 121.6017994420231
$wsRaw.Range("C5").Value = 6.226907061890508
$wsRaw.Range("D5").Value = 6.224037891463015
$wsRaw.Range("E5").Value = 51.55675097770298
$wsRaw.Range("F5").Value = 51.47691576655491
$wsRaw.Range("G5").Value = 5.310947514377122
$wsRaw.Range("H5").Value = 268.8283216842783
$wsRaw.Range("J5").Value = 5.60059003508913
$wsRaw.Range("L5").Value = 2976.848647266396
$wsRaw.Range("M5").Value = 877.0260989985096
$wsRaw.Range("A6").Value = 565
$wsRaw.Range("B6").Value = 168.9319101544706
$wsRaw.Range("C6").Value = 7.404025875251799
$wsRaw.Range("D6").Value = 7.395483172797983
$wsRaw.Range("E6").Value = 61.88301812652696
$wsRaw.Range("F6").Value = 61.79873719887092
$wsRaw.Range("G6").Value = 6.631241009593016
$wsRaw.Range("H6").Value = 348.7742969885265
$wsRaw.Range("J6").Value = 7.26613118726097
$wsRaw.Range("L6").Value = 4891.416440311873
$wsRaw.Range("M6").Value = 1938.046402406422
$wsRaw.Range("A7").Value = 554
$wsRaw.Range("B7").Value = 209.7779859194627
$wsRaw.Range("C7").Value = 8.332232478005091
$wsRaw.Range("D7").Value = 8.314477293407698
$wsRaw.Range("E7").Value = 70.28495482032018
$wsRaw.Range("F7").Value = 70.19377281226988
$wsRaw.Range("G7").Value = 7.709843018029283
$wsRaw.Range("H7").Value = 413.3631685703516
$wsRaw.Range("J7").Value = 8.611732678548993
$wsRaw.Range("L7").Value = 6854.197547012961
$wsRaw.Range("M7").Value = 3131.503440440872
$wsRaw.Range("N7").Value = 530.4238094679079
$wsRaw.Range("B8").Value = 0.5721657080754552
$wsRaw.Range("C8").Value = 0.512114676648408
$wsRaw.Range("D8").Value = 0.5027385963179115
$wsRaw.Range("G8").Value = 0.04559877201070527
$wsRaw.Range("H8").Value = 3.393082218936029
$wsRaw.Range("L8").Value = 74.8942614420387
$wsRaw.Range("A9").Value = 392
$wsRaw.Range("B9").Value = 20.32151034321782
$wsRaw.Range("C9").Value = 3.082983945489418
$wsRaw.Range("D9").Value = 3.082954101582371
$wsRaw.Range("E9").Value = 24.18674884842458
$wsRaw.Range("F9").Value = 24.11484980770635
$wsRaw.Range("G9").Value = 1.576129479872707
$wsRaw.Range("H9").Value = 59.30413071936137
$wsRaw.Range("J9").Value = 1.235502723320029
$wsRaw.Range("L9").Value = 297.7702245715509
$wsRaw.Range("A10").Value = 384
$wsRaw.Range("B10").Value = 48.04698087037325
$wsRaw.Range("C10").Value = 4.789649623770505
$wsRaw.Range("D10").Value = 4.789385732980367
$wsRaw.Range("E10").Value = 38.59153620314448
$wsRaw.Range("F10").Value = 38.51915932782047
$wsRaw.Range("G10").Value = 2.529577127661593
$wsRaw.Range("H10").Value = 117.8199230211342
$wsRaw.Range("J10").Value = 2.454581729606963
$wsRaw.Range("L10").Value = 920.0028696739272
$wsRaw.Range("A11").Value = 375
$wsRaw.Range("B11").Value = 85.27790412697688
$wsRaw.Range("C11").Value = 6.457121087095595
$wsRaw.Range("D11").Value = 6.455638597691814
$wsRaw.Range("E11").Value = 51.73444978279363
$wsRaw.Range("F11").Value = 51.66171090147062
$wsRaw.Range("G11").Value = 3.638795998012128
$wsRaw.Range("H11").Value = 185.8419489477265
$wsRaw.Range("J11").Value = 3.871707269744302
$wsRaw.Range("L11").Value = 2089.500232641966
$wsRaw.Range("M11").Value = 659.9120854924157
$wsRaw.Range("B12").Value = 122.978023673544
$wsRaw.Range("C12").Value = 7.859657619844877
$wsRaw.Range("D12").Value = 7.854838548277741
$wsRaw.Range("E12").Value = 62.19805114652858
$wsRaw.Range("F12").Value = 62.12339795175566
$wsRaw.Range("G12").Value = 4.667939004912062
$wsRaw.Range("H12").Value = 247.9794469768551
$wsRaw.Range("J12").Value = 5.166238478684482
$wsRaw.Range("L12").Value = 3570.29542632243
$wsRaw.Range("M12").Value = 1519.75895644745
$wsRaw.Range("A13").Value = 360
$wsRaw.Range("B13").Value = 159.8118182681903
$wsRaw.Range("C13").Value = 9.021723266318652
$wsRaw.Range("D13").Value = 9.010787872514994
$wsRaw.Range("E13").Value = 70.60581276820122
$wsRaw.Range("F13").Value = 70.52696050898382
$wsRaw.Range("G13").Value = 5.639628513197215
$wsRaw.Range("H13").Value = 305.1702197634342
$wsRaw.Range("J13").Value = 6.357712911738212
$wsRaw.Range("L13").Value = 5234.115318245908
$wsRaw.Range("M13").Value = 2571.819542471451
$wsRaw.Range("N13").Value = 1820.90349058675
$wsRaw.Range("B14").Value = 0.8411720815708631
$wsRaw.Range("C14").Value = 0.5069943332661339
$wsRaw.Range("D14").Value = 0.4968057160504162
$wsRaw.Range("G14").Value = 0.05367122386096896
$wsRaw.Range("H14").Value = 5.008195045323141
$wsRaw.Range("L14").Value = 112.333564880946
$wsRaw.Range("A15").Value = 585
$wsRaw.Range("B15").Value = 21.44808834741062
$wsRaw.Range("C15").Value = 2.592702385502407
$wsRaw.Range("D15").Value = 2.592673664081908
$wsRaw.Range("E15").Value = 24.18375297992464
$wsRaw.Range("F15").Value = 24.10702014466022
$wsRaw.Range("G15").Value = 1.368501376001441
$wsRaw.Range("H15").Value = 67.02377684653456
$wsRaw.Range("J15").Value = 1.367832180541522
$wsRaw.Range("L15").Value = 346.1578730173044
$wsRaw.Range("A16").Value = 575
$wsRaw.Range("B16").Value = 55.69596680848827
$wsRaw.Range("C16").Value = 4.214191400883609
$wsRaw.Range("D16").Value = 4.214070745712717
$wsRaw.Range("E16").Value = 38.85915646213111
$wsRaw.Range("F16").Value = 38.791585719158
$wsRaw.Range("G16").Value = 2.705573500517634
$wsRaw.Range("H16").Value = 143.6594589645436
$wsRaw.Range("J16").Value = 2.931825693153951
$wsRaw.Range("L16").Value = 1097.634637926815
$wsRaw.Range("A17").Value = 561
$wsRaw.Range("B17").Value = 97.54561141435521
$wsRaw.Range("C17").Value = 5.646232302827041
$wsRaw.Range("D17").Value = 5.645292669087226
$wsRaw.Range("E17").Value = 53.28495199827061
$wsRaw.Range("F17").Value = 53.21599421038238
$wsRaw.Range("G17").Value = 3.910568374110508
$wsRaw.Range("H17").Value = 224.1483997736751
$wsRaw.Range("J17").Value = 4.574457138238267
$wsRaw.Range("L17").Value = 2484.184793904862
$wsRaw.Range("M17").Value = 583.8147108601621
$wsRaw.Range("A18").Value = 552
$wsRaw.Range("B18").Value = 139.7959448178068
$wsRaw.Range("C18").Value = 6.814187089268948
$wsRaw.Range("D18").Value = 6.810958289604335
$wsRaw.Range("E18").Value = 67.50094061914709
$wsRaw.Range("F18").Value = 67.42792092452903
$wsRaw.Range("G18").Value = 5.048660016970894
$wsRaw.Range("H18").Value = 298.2418834332651
$wsRaw.Range("J18").Value = 6.08656904965847
$wsRaw.Range("L18").Value = 4423.967842210355
$wsRaw.Range("M18").Value = 1567.601450377289
$wsRaw.Range("A19").Value = 541
$wsRaw.Range("B19").Value = 176.9981058109088
$wsRaw.Range("C19").Value = 7.745009987747919
$wsRaw.Range("D19").Value = 7.737764072804728
$wsRaw.Range("E19").Value = 80.35275832372183
$wsRaw.Range("F19").Value = 80.27318757498355
$wsRaw.Range("G19").Value = 5.999774663278867
$wsRaw.Range("H19").Value = 358.9860285040294
$wsRaw.Range("J19").Value = 7.32624547967407
$wsRaw.Range("L19").Value = 6614.152094010635
$wsRaw.Range("M19").Value = 2811.154536816042

# --- Summary sheet updates ---
$wsSummary.Range("C2").Value = 554
$wsSummary.Range("D2").Value = 6854.197547012961
$wsSummary.Range("E2").Value = 8.314477293407698
$wsSummary.Range("F2").Value = 70.19377281226988
$wsSummary.Range("C3").Value = 360
$wsSummary.Range("D3").Value = 5234.115318245908
$wsSummary.Range("E3").Value = 9.010787872514994
$wsSummary.Range("F3").Value = 70.52696050898382
$wsSummary.Range("C4").Value = 541
$wsSummary.Range("D4").Value = 6614.152094010635
$wsSummary.Range("E4").Value = 7.737764072804728
$wsSummary.Range("F4").Value = 80.27318757498355
